# Append three new observation rows (43-45) to the "Artfynd" export sheet,
# matching the source data feed's most recent rows.
#
# Column layout (row 1 header): A=Id, B=Taxonsorteringsordning, C=Valideringsstatus,
# D=Rödlistade, E=TaxonId, F=Artnamn, G=Vetenskapligt namn, H=Auktor, I=Antal,
# J=Enhet, K=Ålder-Stadium, L=Kön, M=Aktivitet, N=Metod, O=Huvudlokal,
# P=Lokalnamn, Q=Ost, R=Nord, S=Noggrannhet, T=Län, U=Kommun, V=Provins,
# W=Socken, X=Externid, Y=Startdatum, Z=Starttid, AA=Slutdatum, AB=Sluttid,
# AC=Publik kommentar, AD=Ej återfunnen, AE=Osäker artbestämning,
# AF=Bestämningsmetod, AG=Ospontan, ... AT=Bestämningsår, AW=Rapportör,
# AX=Observatörer, AY=Projektnamn.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe forces Excel to keep a numeric/date-looking literal
# ("1", "2026-02-14", ...) as text instead of silently parsing it into a
# number/date serial -- this mirrors the source file, where every value in
# these columns is stored as plain text (inline string), never as a number
# or date. A bare apostrophe likewise yields an explicit empty-text cell
# (rather than leaving the cell completely unset) for the columns whose
# source rows carry an (empty) text placeholder.

function Set-Row($r, $vals) {
    foreach ($col in $vals.Keys) {
        $ws.Cells.Item($r, $col).Value = $vals[$col]
    }
}

# Row 43: Tjäder (Tetrao urogallus), 2026-02-14, Sims bodarna
Set-Row 43 @{
    1  = 131219323          # A  Id
    2  = 57073               # B  Taxonsorteringsordning
    4  = "LC"                 # D  Rödlistade
    5  = 100138               # E  TaxonId
    6  = "Tjäder"             # F  Artnamn
    7  = "Tetrao urogallus"   # G  Vetenskapligt namn
    8  = "Linnaeus, 1758"     # H  Auktor
    9  = "'1"                 # I  Antal (text "1")
    11 = "'"                  # K  Ålder-Stadium (empty text)
    12 = "'"                  # L  Kön (empty text)
    13 = "förbiflygande"      # M  Aktivitet
    14 = "'"                  # N  Metod (empty text)
    16 = "Sims bodarna, Dlr"  # P  Lokalnamn
    17 = 515391                # Q  Ost
    18 = 6704950                # R  Nord
    19 = 25                   # S  Noggrannhet
    20 = "Dalarna"            # T  Län
    21 = "Borlänge"           # U  Kommun
    22 = "Dalarna"            # V  Provins
    23 = "Stora Tuna"         # W  Socken
    25 = "'2026-02-14"        # Y  Startdatum
    27 = "'2026-02-14"        # AA Slutdatum
    30 = $false                # AD Ej återfunnen
    31 = $false                # AE Osäker artbestämning
    33 = $false                # AG Ospontan
    46 = "'"                  # AT Bestämningsår (empty text)
    49 = "Anna-Lena Thommson" # AW Rapportör
    50 = "Anna-Lena Thommson" # AX Observatörer
    51 = "'"                  # AY Projektnamn (empty text)
}

# Row 44: Tretåig hackspett (Picoides tridactylus), 2026-02-14, Sims bodarna
Set-Row 44 @{
    1  = 131219238
    2  = 57884
    4  = "NT"
    5  = 100109
    6  = "Tretåig hackspett"
    7  = "Picoides tridactylus"
    8  = "(Linnaeus, 1758)"
    9  = "'"                  # I  Antal (empty text)
    11 = "'"                  # K  Ålder-Stadium (empty text)
    12 = "'"                  # L  Kön (empty text)
    13 = "färska spår"        # M  Aktivitet
    14 = "'"                  # N  Metod (empty text)
    16 = "Sims bodarna, Dlr"
    17 = 515365
    18 = 6704964
    19 = 50
    20 = "Dalarna"
    21 = "Borlänge"
    22 = "Dalarna"
    23 = "Stora Tuna"
    25 = "'2026-02-14"
    27 = "'2026-02-14"
    29 = "Ringhack på tall."  # AC Publik kommentar
    30 = $false
    31 = $false
    33 = $false
    46 = "'"
    49 = "Anna-Lena Thommson"
    50 = "Anna-Lena Thommson"
    51 = "'"
}

# Row 45: Tretåig hackspett (Picoides tridactylus), 2026-02-14, Sims bodarna
Set-Row 45 @{
    1  = 131219286
    2  = 57884
    4  = "NT"
    5  = 100109
    6  = "Tretåig hackspett"
    7  = "Picoides tridactylus"
    8  = "(Linnaeus, 1758)"
    9  = "'"
    11 = "'"
    12 = "'"
    13 = "färska spår"
    14 = "'"
    16 = "Sims bodarna, Dlr"
    17 = 515382
    18 = 6704962
    19 = 50
    20 = "Dalarna"
    21 = "Borlänge"
    22 = "Dalarna"
    23 = "Stora Tuna"
    25 = "'2026-02-14"
    27 = "'2026-02-14"
    29 = "Ringhack på tall."
    30 = $false
    31 = $false
    33 = $false
    46 = "'"
    49 = "Anna-Lena Thommson"
    50 = "Anna-Lena Thommson"
    51 = "'"
}
